$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Unknown Title"
$ws.Range("D2").Value = "Unknown Abstract"
$ws.Range("E2").Value = "[]"
$ws.Range("F2").Value = "not found"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = "1970-01-01"
$ws.Range("E3").Value = "[Tao%Chen%NULL%0, Di%Wu%NULL%0, Huilong%Chen%NULL%0, Weiming%Yan%NULL%0, Danlei%Yang%NULL%0, Guang%Chen%NULL%0, Ke%Ma%NULL%0, Dong%Xu%NULL%0, Haijing%Yu%NULL%0, Hongwu%Wang%NULL%0, Tao%Wang%NULL%0, Wei%Guo%NULL%0, Jia%Chen%NULL%0, Chen%Ding%NULL%0, Xiaoping%Zhang%NULL%0, Jiaquan%Huang%NULL%0, Meifang%Han%NULL%0, Shusheng%Li%NULL%0, Xiaoping%Luo%NULL%0, Jianping%Zhao%NULL%0, Qin%Ning%NULL%0]"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "BMJ Publishing Group Ltd."
$ws.Range("D4").Value = "Background`nid=""Par1"">A COVID-19 outbreak started in Wuhan, China, last December and now has become a global pandemic.`n`n The clinical information in caring of critically ill patients with COVID-19 needs to be shared timely, especially under the situations that there is still a largely ongoing spread of COVID-19 in many countries.`n`n`nMethods`nid=""Par2"">A multicenter prospective observational study investigated all the COVID-19 patients received in 19 ICUs of 16 hospitals in Wuhan, China, over 24 h between 8 AM February 2h and 8 AM February 27, 2020. The demographic information, clinical characteristics, vital signs, complications, laboratory values, and clinical managements of the patients were studied.`n`n`nResults`nid=""Par3"">A total of 226 patients were included.`n`n Their median (interquartile range, IQR) age was 64 (57–70) years, and 139 (61.5%) patients were male.`n`n The duration from the date of ICU admission to the study date was 11 (5–17) days, and the duration from onset of symptoms to the study date was 31 (24–36) days.`n`n Among all the patients, 155 (68.6%) had at least one coexisting disease, and their sequential organ failure assessment score was 4 (2–8).`n`n Organ function damages were found in most of the patients: ARDS in 161 (71.2%) patients, septic shock in 34 (15.0%) patients, acute kidney injury occurred in 57 (25.2%) patients, cardiac injury in 61 (27.0%) patients, and lymphocytopenia in 160 (70.8%) patients.`n`n Of all the studied patients, 85 (37.6%) received invasive mechanical ventilation, including 14 (6.2%) treated with extracorporeal membrane oxygenation (ECMO) at the same time, 20 (8.8%) received noninvasive mechanical ventilation, and 24 (10.6%) received continuous renal replacement therapy.`n`n By April 9, 2020, 87 (38.5%) patients were deceased and 15 (6.7%) were still in the hospital.`n`n`nConclusions`nid=""Par4"">Critically ill patients with COVID-19 are associated with a higher risk of severe complications and need to receive an intensive level of treatments.`n`n COVID-19 poses a great strain on critical care resources in hospitals.`n`n`nTrial registration`nid=""Par5"">Chinese Clinical Trial Registry, ChiCTR2000030164. Registered on February 24, 2020, http://www.`n`nchictr.`n`norg.`n`ncn/edit.`n`naspx?pid=49983&amp;htm=4`n"
$ws.Range("E4").Value = "[Yuan%Yu%NULL%0, Dan%Xu%NULL%0, Shouzhi%Fu%NULL%0, Jun%Zhang%NULL%0, Xiaobo%Yang%NULL%0, Liang%Xu%NULL%0, Jiqian%Xu%NULL%0, Yongran%Wu%NULL%0, Chaolin%Huang%NULL%0, Yaqi%Ouyang%NULL%0, Luyu%Yang%NULL%0, Minghao%Fang%NULL%0, Hongwen%Xiao%NULL%0, Jing%Ma%NULL%0, Wei%Zhu%NULL%0, Song%Hu%NULL%0, Quan%Hu%NULL%0, Daoyin%Ding%NULL%0, Ming%Hu%NULL%0, Guochao%Zhu%NULL%0, Weijiang%Xu%NULL%0, Jun%Guo%NULL%0, Jinglong%Xu%NULL%0, Haitao%Yuan%NULL%0, Bin%Zhang%NULL%0, Zhui%Yu%yuzhui@whu.edu.cn%0, Dechang%Chen%icudechangchen@163.com%0, Shiying%Yuan%yuan_shiying@163.com%0, You%Shang%you_shanghust@163.com%0]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "BioMed Central"
$ws.Range("C5").Value = "Unknown Title"
$ws.Range("D5").Value = "Unknown Abstract"
$ws.Range("E5").Value = "[]"
$ws.Range("F5").Value = "not found"
$ws.Range("G5").Value = "N/A"
$ws.Range("I5").Value = ""
$ws.Range("C6").Value = "Unknown Title"
$ws.Range("D6").Value = "Unknown Abstract"
$ws.Range("E6").Value = "[]"
$ws.Range("F6").Value = "not found"
$ws.Range("G6").Value = "N/A"
$ws.Range("I6").Value = ""
$ws.Range("E7").Value = "[Pavan K.%Bhatraju%NULL%0, Bijan J.%Ghassemieh%NULL%0, Michelle%Nichols%NULL%0, Richard%Kim%NULL%0, Keith R.%Jerome%NULL%0, Arun K.%Nalla%NULL%0, Alexander L.%Greninger%NULL%0, Sudhakar%Pipavath%NULL%0, Mark M.%Wurfel%NULL%0, Laura%Evans%NULL%0, Patricia A.%Kritek%NULL%0, T. Eoin%West%NULL%0, Andrew%Luks%NULL%0, Anthony%Gerbino%NULL%0, Chris R.%Dale%NULL%0, Jason D.%Goldman%NULL%0, Shane%O’Mahony%NULL%0, Carmen%Mikacenic%NULL%0]"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = "Massachusetts Medical Society"
$ws.Range("D8").Value = "We investigated the clinical course of individuals with 2019 novel coronavirus disease (COVID-19) who were transferred from the Diamond Princess cruise ship to 12 local hospitals.`n The conditions and clinical courses of patients with pneumonia were compared with those of patients without pneumonia.`n Among 70 patients (median age: 67 years) analyzed, the major symptoms were fever (64.3%), cough (54.3%), and general fatigue (24.3%).`n Forty-three patients (61.4%) had pneumonia.`n Higher body temperature, heart rate, and respiratory rate as well as higher of lactate dehydrogenase (LDH), aspartate aminotransferase (AST), and C-reactive protein (CRP) levels and lower serum albumin level and lymphocyte count were associated with the presence of pneumonia.`n Ground-glass opacity was found in 97.7% of the patients with pneumonia.`n Patients were administered neuraminidase inhibitors (20%), lopinavir/ritonavir (32.9%), and ciclesonide inhalation (11.4%).`n Mechanical ventilation and veno-venous extracorporeal membrane oxygenation was performed on 14 (20%) and 2 (2.9%) patients, respectively; two patients died.`n The median duration of intubation was 12 days.`n The patients with COVID-19 transferred to local hospitals during the outbreak had severe conditions and needed close monitoring.`n The severity of COVID-19 depends on the presence of pneumonia.`n High serum LDH, AST and CRP levels and low serum albumin level and lymphocyte count were found to be predictors of pneumonia.`n It was challenging for local hospitals to admit and treat these patients during the outbreak of COVID-19. Assessment of severity was crucial to manage a large number of patients.`n"
$ws.Range("E8").Value = "[Hideaki%Kato%NULL%0, Hiroyuki%Shimizu%NULL%1, Yasushi%Shibue%NULL%1, Tomohiro%Hosoda%NULL%1, Keisuke%Iwabuchi%NULL%1, Kotaro%Nagamine%NULL%1, Hiroki%Saito%NULL%1, Reimin%Sawada%NULL%1, Takayuki%Oishi%NULL%1, Jun%Tsukiji%NULL%1, Hiroyuki%Fujita%NULL%1, Ryosuke%Furuya%NULL%1, Makoto%Masuda%NULL%1, Osamu%Akasaka%NULL%1, Yu%Ikeda%NULL%1, Mitsuo%Sakamoto%NULL%1, Kazuya%Sakai%NULL%1, Munehito%Uchiyama%NULL%1, Hiroki%Watanabe%NULL%1, Nobuhiro%Yamaguchi%NULL%1, Ryoko%Higa%NULL%1, Akiko%Sasaki%NULL%1, Katsuaki%Tanaka%NULL%1, Yukitoshi%Toyoda%NULL%1, Shinsuke%Hamanaka%NULL%1, Naoki%Miyazawa%NULL%1, Atsuko%Shimizu%NULL%1, Fumie%Fukase%NULL%1, Shunsuke%Iwai%NULL%1, Yuko%Komase%NULL%1, Tsutomu%Kawasaki%NULL%1, Isao%Nagata%NULL%1, Yusuke%Nakayama%NULL%1, Tetsuhiro%Takei%NULL%1, Katsuo%Kimura%NULL%1, Reiko%Kunisaki%NULL%1, Makoto%Kudo%NULL%1, Ichiro%Takeuchi%NULL%0, Hideaki%Nakajima%NULL%1]"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "Japanese Society of Chemotherapy and The Japanese Association for Infectious Diseases. Published by Elsevier Ltd."
$ws.Range("C9").Value = "Unknown Title"
$ws.Range("E9").Value = "[]"
$ws.Range("F9").Value = "not found"
$ws.Range("G9").Value = "N/A"
$ws.Range("H9").Value = "1970-01-01"
$ws.Range("J9").Value = ""
$ws.Range("E10").Value = "[Xiaobo%Yang%NULL%0, Yuan%Yu%NULL%0, Jiqian%Xu%NULL%0, Huaqing%Shu%NULL%0, Jia'an%Xia%NULL%0, Hong%Liu%NULL%0, Yongran%Wu%NULL%0, Lu%Zhang%NULL%0, Zhui%Yu%NULL%0, Minghao%Fang%NULL%0, Ting%Yu%NULL%0, Yaxin%Wang%NULL%0, Shangwen%Pan%NULL%0, Xiaojing%Zou%NULL%0, Shiying%Yuan%NULL%0, You%Shang%NULL%0]"
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = "Elsevier Ltd."
$ws.Range("D11").Value = "Background: The risk factors for adverse events of Coronavirus Disease-19 (COVID-19) have not been well described.`n We aimed to explore the predictive value of clinical, laboratory and CT imaging characteristics on admission for short-term outcomes of COVID-19 patients.`n"
$ws.Range("E11").Value = "[Peng Peng%Xu%NULL%0, Rong Hua%Tian%NULL%1, Song%Luo%NULL%0, Zi Yue%Zu%NULL%1, Bin%Fan%NULL%1, Xi Ming%Wang%NULL%1, Kai%Xu%NULL%1, Jiang Tao%Wang%NULL%1, Juan%Zhu%NULL%1, Ji Chan%Shi%NULL%1, Feng%Chen%NULL%1, Bing%Wan%NULL%1, Zhi Han%Yan%NULL%1, Rong Pin%Wang%NULL%1, Wen%Chen%NULL%1, Wen Hui%Fan%NULL%1, Can%Zhang%NULL%1, Meng Jie%Lu%NULL%1, Zhi Yuan%Sun%NULL%0, Chang Sheng%Zhou%NULL%0, Li Na%Zhang%NULL%1, Fei%Xia%NULL%0, Li%Qi%NULL%1, Wei%Zhang%NULL%0, Jing%Zhong%NULL%1, Xiao Xue%Liu%NULL%1, Qi Rui%Zhang%NULL%1, Guang Ming%Lu%NULL%0, Long Jiang%Zhang%NULL%0]"
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = "Ivyspring International Publisher"
$ws.Range("E12").Value = "[Matthew J%Cummings%NULL%0, Matthew R%Baldwin%NULL%0, Darryl%Abrams%NULL%2, Samuel D%Jacobson%NULL%2, Benjamin J%Meyer%NULL%2, Elizabeth M%Balough%NULL%2, Justin G%Aaron%NULL%2, Jan%Claassen%NULL%0, LeRoy E%Rabbani%NULL%2, Jonathan%Hastie%NULL%2, Beth R%Hochman%NULL%2, John%Salazar-Schicchi%NULL%2, Natalie H%Yip%NULL%2, Daniel%Brodie%NULL%0, Max R%O'Donnell%NULL%2]"
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = "Elsevier Ltd."
$ws.Range("E13").Value = "[Michael G%Argenziano%NULL%0, Samuel L%Bruce%NULL%0, Cody L%Slater%NULL%0, Jonathan R%Tiao%NULL%0, Matthew R%Baldwin%NULL%0, R Graham%Barr%NULL%0, Bernard P%Chang%NULL%0, Katherine H%Chau%NULL%0, Justin J%Choi%NULL%0, Nicholas%Gavin%NULL%0, Parag%Goyal%NULL%0, Angela M%Mills%NULL%0, Ashmi A%Patel%NULL%0, Marie-Laure S%Romney%NULL%0, Monika M%Safford%NULL%0, Neil W%Schluger%NULL%0, Soumitra%Sengupta%NULL%0, Magdalena E%Sobieszczyk%NULL%0, Jason E%Zucker%NULL%0, Paul A%Asadourian%NULL%0, Fletcher M%Bell%NULL%0, Rebekah%Boyd%NULL%0, Matthew F%Cohen%NULL%0, MacAlistair I%Colquhoun%NULL%0, Lucy A%Colville%NULL%0, Joseph H%de Jonge%NULL%0, Lyle B%Dershowitz%NULL%0, Shirin A%Dey%NULL%0, Katherine A%Eiseman%NULL%0, Zachary P%Girvin%NULL%0, Daniella T%Goni%NULL%0, Amro A%Harb%NULL%0, Nicholas%Herzik%NULL%0, Sarah%Householder%NULL%0, Lara E%Karaaslan%NULL%0, Heather%Lee%NULL%0, Evan%Lieberman%NULL%0, Andrew%Ling%NULL%0, Ree%Lu%NULL%0, Arthur Y%Shou%NULL%0, Alexander C%Sisti%NULL%0, Zachary E%Snow%NULL%0, Colin P%Sperring%NULL%0, Yuqing%Xiong%NULL%0, Henry W%Zhou%NULL%0, Karthik%Natarajan%NULL%0, George%Hripcsak%NULL%0, Ruijun%Chen%NULL%0]"
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = "BMJ Publishing Group Ltd."
$ws.Range("C14").Value = "Unknown Title"
$ws.Range("D14").Value = "Unknown Abstract"
$ws.Range("E14").Value = "[]"
$ws.Range("F14").Value = "not found"
$ws.Range("G14").Value = "N/A"
$ws.Range("H14").Value = "1970-01-01"
$ws.Range("I14").Value = ""
$ws.Range("E15").Value = "[Sara C.%Auld%NULL%0, Mark%Caridi-Scheible%NULL%0, James M.%Blum%NULL%0, Chad%Robichaux%NULL%0, Colleen%Kraft%NULL%0, Jesse T.%Jacob%NULL%0, Craig S.%Jabaley%NULL%0, David%Carpenter%NULL%0, Roberta%Kaplow%NULL%0, Alfonso C.%Hernandez-Romieu%NULL%0, Max W.%Adelman%NULL%0, Greg S.%Martin%NULL%0, Craig M.%Coopersmith%NULL%0, David J.%Murphy%NULL%0, NULL%NULL%NULL%0, NULL%NULL%NULL%0]"
$ws.Range("I15").Value = ""
$ws.Range("J15").Value = "Lippincott Williams &amp; Wilkins"
$ws.Range("E16").Value = "[Kevin%Hur%NULL%0, Caroline P. E.%Price%NULL%1, Elizabeth L.%Gray%NULL%1, Reeti K.%Gulati%NULL%1, Matthew%Maksimoski%NULL%1, Samuel D.%Racette%NULL%1, Alexander L.%Schneider%NULL%1, Ashoke R.%Khanwalkar%NULL%1]"
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = "SAGE Publications"
